$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("averagetransport work")
$ws2 = $wb.Worksheets.Item("shipfuel eff")

# ---------------------------------------------------------------------------
# Sheet 1: "averagetransport work"
# ---------------------------------------------------------------------------

# New conversion-factor block (I2:J3)
$ws1.Range("I2").Value = 1
$ws1.Range("J2").Value = "NM"
$ws1.Range("I3").Value = 1.852
$ws1.Range("J3").Value = "km"

# Unit label changed from "tonne-NM" to "gigatonne-NM"
$ws1.Range("B4").Value = "gigatonne-NM"

# Row 6: values rescaled from raw tonne-NM counts down to gigatonne-NM, and
# the old scientific-notation formatting is cleared.
$ws1.Range("B6:F6").ClearFormats()
$ws1.Range("B6").Value = 3
$ws1.Range("C6").Value = 7
$ws1.Range("D6").Value = 0.5
$ws1.Range("E6").Value = 4
$ws1.Range("F6").Value = 0.3

# Row 7: new derived row converting to km using the factor in I3, formatted
# to one decimal place.
$ws1.Range("B7:F7").NumberFormat = "0.0"
$ws1.Range("B7").Formula = "=+B6*`$I`$3"
$ws1.Range("C7").Formula = "=+C6*`$I`$3"
$ws1.Range("D7").Formula = "=+D6*`$I`$3"
$ws1.Range("E7").Formula = "=+E6*`$I`$3"
$ws1.Range("F7").Formula = "=+F6*`$I`$3"

# Rows 45-46: a new 35-column "ship efficiency" table duplicated from sheet 2.
$fuelHeaders = @("T_MFO","B_MFO","G_MFO","C_MFO","O_MFO","T_SCR","B_SCR","G_SCR","C_SCR","O_SCR","T_LNG","B_LNG","G_LNG","C_LNG","O_LNG","T_ELC","B_ELC","G_ELC","C_ELC","O_ELC","T_MET","B_MET","G_MET","C_MET","O_MET","T_LPG","B_LPG","G_LPG","C_LPG","O_LPG","T_HYD","B_HYD","G_HYD","C_HYD","O_HYD")
$blockVals = @(5.556, 12.964, 0.92600000000000005, 7.4080000000000004, 0.55559999999999998)

$col = 2
foreach ($h in $fuelHeaders) {
    $ws1.Cells.Item(45, $col).ClearFormats()
    $ws1.Cells.Item(45, $col).Value = $h
    $col++
}

$ws1.Range("B46:AJ46").NumberFormat = "0.0"
$col = 2
for ($grp = 0; $grp -lt 7; $grp++) {
    foreach ($v in $blockVals) {
        $ws1.Cells.Item(46, $col).Value = $v
        $col++
    }
}

# ---------------------------------------------------------------------------
# Sheet 2: "shipfuel eff"
# ---------------------------------------------------------------------------

# New units label above the header row
$ws2.Range("C4").Value = "ton*km/MJ"

# Real (non-placeholder) data replacing the old 0.1 / 0.2 stand-in values.
$ws2.Range("H6").Value = 10
$ws2.Range("I6").Value = 10
$ws2.Range("J6").Value = 5
$ws2.Range("K6").Value = 5
$ws2.Range("L6").Value = 5

$ws2.Range("C7").Value = 10
$ws2.Range("D7").Value = 10
$ws2.Range("E7").Value = 5
$ws2.Range("F7").Value = 5
$ws2.Range("G7").Value = 5

$ws2.Range("M8").Value = 10
$ws2.Range("N8").Value = 10
$ws2.Range("O8").Value = 5
$ws2.Range("P8").Value = 5
$ws2.Range("Q8").Value = 5

$ws2.Range("AB9").Value = 10
$ws2.Range("AC9").Value = 10
$ws2.Range("AD9").Value = 5
$ws2.Range("AE9").Value = 5
$ws2.Range("AF9").Value = 5

$ws2.Range("W10").Value = 10
$ws2.Range("X10").Value = 10
$ws2.Range("Y10").Value = 5
$ws2.Range("Z10").Value = 5
$ws2.Range("AA10").Value = 5

$ws2.Range("AG11").Value = 10
$ws2.Range("AH11").Value = 10
$ws2.Range("AI11").Value = 5
$ws2.Range("AJ11").Value = 5
$ws2.Range("AK11").Value = 5

$ws2.Range("AG12").Value = 10
$ws2.Range("AH12").Value = 10
$ws2.Range("AI12").Value = 5
$ws2.Range("AJ12").Value = 5
$ws2.Range("AK12").Value = 5

$ws2.Range("R14").Value = 10
$ws2.Range("S14").Value = 10
$ws2.Range("T14").Value = 5
$ws2.Range("U14").Value = 5
$ws2.Range("V14").Value = 5

# ---------------------------------------------------------------------------
# View state: sheet 1 becomes the active / selected tab, sheet 2 keeps a
# plain (non-active) selection.
# ---------------------------------------------------------------------------
$ws2.Range("C5").Select()
$ws1.Activate()
$ws1.Range("F42").Select()
